# Apply "add main BoM and images" changes to the RP2040-VCO BoM workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BoM")

# 1. References list for the 10k resistor row gains R17.
$ws.Range("C28").Value = "R8 R10 R12 R13 R14 R17"

# 2. Quantity Per PCB for that same row goes from 5 to 6 (R17 added).
#    The BoM stores this column as text (shared string), like every other
#    row, so copy an existing "6" text cell over via PasteSpecial(values)
#    rather than assigning a literal which Excel would coerce to a number.
$ws.Range("A14").Copy()
$ws.Range("F28").PasteSpecial(-4163)

# 3. Generation timestamp bumped for the regenerated BoM.
$ws.Range("B5").Value = "2023-04-09_19-20-43"

# 4. Component Count / Fitted Components summary text updated
#    (both cells share the same underlying text in the original file).
$ws.Range("D3").Value = "53 (39 SMD/ 14 THT)"
$ws.Range("D4").Value = "53 (39 SMD/ 14 THT)"

# 5. Total Components count increments by one.
$ws.Range("D6").Value = 53

# 6. Widen the References column to fit the longer reference list
#    (target OOXML width 27.7109375; ColumnWidth rounds to the nearest
#    whole pixel on write, so 26.83 is the input that lands closest).
$ws.Columns.Item(3).ColumnWidth = 26.83
